$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two more benchmark rows were added under "Consolidate spans" ("Consolidate
# spans2" / "Consolidate spans3"), continuing the existing FPS table and its
# shared C (=B/30) and D (=B/$B$2) formulas.
$ws.Range("A30").Value = "Consolidate spans2"
$ws.Range("B30").Value = 343

$ws.Range("A31").Value = "Consolidate spans3"
$ws.Range("B31").Value = 354

$ws.Range("C30:C31").Formula = "=B30/30"
$ws.Range("D30:D31").Formula = "=B30/`$B`$2"

# Move the view: the user had scrolled/selected further down the sheet.
$ws.Range("A31").Select()

$wb.Save()
